$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T230004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2007%3A00%3A04%20GMT&X-Amz-Signature=b070701ebb6800e5ed1a3271c0165040c3a14f0449736f0e7555b8ef57ecc8f3'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T230004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2007%3A00%3A04%20GMT&X-Amz-Signature=c48a182891c1067858aa2231f3322fadb6f888034dc24c33bf9ed7323c9dccee'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T230004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2007%3A00%3A04%20GMT&X-Amz-Signature=3819e1cebd3cc67dd8ab80a4b906a1856fab0542b231fc195ecbe1d49caa41d5'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T230004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2007%3A00%3A04%20GMT&X-Amz-Signature=f2f488e2a95e8b46ec061e7ca168ca074d53c31a87811e2e853cd7f51534b77d'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSHRENQIR2%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T230004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC4aCXVzLXdlc3QtMiJGMEQCIGf4PXk5mIZnHcaLIp53%2FRViRLEN4IiwYlcVoDZztYrpAiAS6Gep2RJ0cK3j82B6ED4z%2FjbXkWjg1TGrmYvfBB94EirmAwiX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMT0L43uAwaH2mrmpIKroDdx0a5ypwL0KfzEQoAtOfX%2FRtXUzxZTwJFTgGw9ICNKAcK1Ua7rMWkOV4WjazXW7v897wqOmDc0X0YN0lPULWX9VM3heRWs98hhNRfR2vTCNjydQhgV4UwtPcFFyE%2BfzGl5fVJkoiNxF5iD4jflTxqF7hwYnTaXqOhAouYaYjR62CWOcCaHMokWAU8Ik%2FR387t3lcwRKPzWTy6wBZGEcoDLumAFaaGj69LfzgxkqpQXzNVHv6HKOr5C3Uw5kseFwnnY6mHxYoxZlb7lXR4qLsp3z%2BL9mBptZ4AlotxAEoqx1jiZcvi91hGC1yoWmCjOi1PoV9cysQQruK3O49dFKlcvemawNg%2FAZt304oJpvugesrkvbWg1Xu3e8mB1z4xgQIJ2oBoWGlkxZ73ifKvIwutjUDNCz0aww6IxoDfJJR82MNZqbk5hjXa7xGkZE%2BOeuefxc%2Bnr%2B%2BczdOnIxPpLT1T8AkVBIHdUGj%2BVxdA0IB8LSuohpJITSVYaMBvl96o%2F%2BRu77VmUNOrKM3Ri0zSTJknBdxPk%2FFOpxZ%2FEiaO40Ie3Jwl%2FbQEByY0LWkL0BjFBqWw5tWEQ76pOHhgTDp76a%2FBjqmAUlYTKltUEZA7KmIBJUyujqtGZTfNL9hEH9HkjwzUV0Dn9wQZjRZp8r75CwdBYViV%2BxD0FijFhrk6vkYNob9pDH2NtYhdIXGJQxLZJ7CC4jzPAq1%2BEEdhYYp3kZ3pmhLq0u4%2B4L93HdDTxqV7lvyBmRxDWLl7sSbNOVl7CzVXnN8dSvmJxG%2BEREFIG6yTGo94rW1bI6XaiBTvks4BjP9BLIROsvYSMY%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2007%3A00%3A04%20GMT&X-Amz-Signature=14e5efbbb656782fa74ff6cdfd29af4d4f05c05e7748d0c00d7ddd19d3955a7f'
